$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.796.49'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '1.976.69'
$ws.Range("E3").Value = '  +0.80%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.11'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("E6").Value = '  +1.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.10'
$ws.Range("E7").Value = '  +3.67%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +2.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0799'
$ws.Range("E10").Value = '  -0.89%  '

$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.42'
$ws.Range("E12").Value = '  +5.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.847'
$ws.Range("E13").Value = '  +2.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.03'
$ws.Range("E14").Value = '  -0.87%  '

$ws.Range("D15").Value = '2.267.39'
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("D17").Value = '1.977.59'
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("D18").Value = '36.688.97'
$ws.Range("E18").Value = '  +0.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.12'
$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("D20").Value = '0.0₃0860'
$ws.Range("E20").Value = '  -0.14%  '

$ws.Range("E21").Value = '  +1.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.12'
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("E24").Value = '  +1.90%  '

$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("E26").Value = '  +3.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.24'
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.83'
$ws.Range("E28").Value = '  +1.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.45'
$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.36'
$ws.Range("E30").Value = '  +20.82%  '

$ws.Range("E31").Value = '  +1.87%  '

$ws.Range("E32").Value = '  +2.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0621'
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.50'
$ws.Range("E34").Value = '  +5.62%  '

$ws.Range("E35").Value = '  +1.55%  '

$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.34'
$ws.Range("E37").Value = '  -1.21%  '

$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.45'
$ws.Range("E39").Value = '  -10.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0972'
$ws.Range("E40").Value = '  -3.38%  '

$ws.Range("E41").Value = '  +0.75%  '

$ws.Range("E42").Value = '  +0.90%  '

$ws.Range("E43").Value = '  -0.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.07'
$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = '1.371.91'
$ws.Range("E45").Value = '  +1.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.57'
$ws.Range("E46").Value = '  +2.24%  '

$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("E48").Value = '  +1.20%  '

$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.29'
$ws.Range("E50").Value = '  +6.22%  '

$ws.Range("D51").Value = '2.161.11'
$ws.Range("E51").Value = '  +1.02%  '
